$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the row above (row 33, same banding) down to row 35
# so the new entry keeps the workbook's existing fill/number-format styling.
$ws.Range("C33:F33").Copy()
$ws.Range("C35:F35").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("G4").Copy()
$ws.Range("G35").PasteSpecial(-4122) # xlPasteFormats
$ws.Application.CutCopyMode = 0

# Fill in row 35 (entry #33) with the new expense record
$ws.Range("C35").Value = "支出"
$ws.Range("D35").Value = 300
$ws.Range("E35").Value2 = 43201
$ws.Range("F35").Value = "生活费"
$ws.Range("G35").Value = "生活费(4/11-4/20)"

# Update the view's scroll position / active selection as recorded in the file
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("E36").Select()
